# Auto-generated edit script applying numeric updates to Exodus_Profits price/profit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (H, I, J, K, L, M, N columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 970.9259
$ws.Range("J17").Value2 = 1079.1666
$ws.Range("L17").Value2 = 3237.4998
$ws.Range("N17").Value2 = -3573.4998
$ws.Range("H33").Value2 = 145.94737
$ws.Range("I33").Value2 = 142.94444
$ws.Range("K33").Value2 = 142.94444
$ws.Range("M33").Value2 = 86.05556000000001
$ws.Range("H80").Value2 = 2047.8334
$ws.Range("I80").Value2 = 415
$ws.Range("J80").Value2 = 3214.1428
$ws.Range("K80").Value2 = 1245
$ws.Range("L80").Value2 = 9642.428400000001
$ws.Range("M80").Value2 = -247
$ws.Range("N80").Value2 = -11638.4284
$ws.Range("H83").Value2 = 2047.8334
$ws.Range("I83").Value2 = 415
$ws.Range("J83").Value2 = 3214.1428
$ws.Range("K83").Value2 = 3735
$ws.Range("L83").Value2 = 28927.2852
$ws.Range("M83").Value2 = 1257
$ws.Range("N83").Value2 = -38911.2852
$ws.Range("H86").Value2 = 3507
$ws.Range("I86").Value2 = 2851
$ws.Range("K86").Value2 = 2851
$ws.Range("M86").Value2 = -1728
$ws.Range("H88").Value2 = 911362.5
$ws.Range("I88").Value2 = 2252349.2
$ws.Range("J88").Value2 = 106770.4
$ws.Range("K88").Value2 = 2252349.2
$ws.Range("L88").Value2 = 106770.4
$ws.Range("M88").Value2 = -2251943.2
$ws.Range("N88").Value2 = -107582.4
$ws.Range("H89").Value2 = 3507
$ws.Range("I89").Value2 = 2851
$ws.Range("K89").Value2 = 14255
$ws.Range("M89").Value2 = -8639
$ws.Range("H91").Value2 = 911362.5
$ws.Range("I91").Value2 = 2252349.2
$ws.Range("J91").Value2 = 106770.4
$ws.Range("K91").Value2 = 2252349.2
$ws.Range("L91").Value2 = 106770.4
$ws.Range("M91").Value2 = -2250945.2
$ws.Range("N91").Value2 = -109578.4
$ws.Range("H96").Value2 = 459.8
$ws.Range("I96").Value2 = 449.75
$ws.Range("J96").Value2 = 500
$ws.Range("K96").Value2 = 1349.25
$ws.Range("L96").Value2 = 1500
$ws.Range("M96").Value2 = 23.75
$ws.Range("N96").Value2 = -4246
$ws.Range("H99").Value2 = 398
$ws.Range("I99").Value2 = 247.5
$ws.Range("J99").Value2 = 1000
$ws.Range("K99").Value2 = 742.5
$ws.Range("L99").Value2 = 3000
$ws.Range("M99").Value2 = 755.5
$ws.Range("N99").Value2 = -5996
$ws.Range("H113").Value2 = 4199.385
$ws.Range("I113").Value2 = 2799.2
$ws.Range("K113").Value2 = 2799.2
$ws.Range("M113").Value2 = 454.8000000000002
$ws.Range("H132").Value2 = 1981.3334
$ws.Range("I132").Value2 = 2122.6667
$ws.Range("J132").Value2 = 1133.3334
$ws.Range("K132").Value2 = 6368.000100000001
$ws.Range("L132").Value2 = 3400.0002
$ws.Range("M132").Value2 = -3838.000100000001
$ws.Range("N132").Value2 = -8460.0002
$ws.Range("H137").Value2 = 485877.47
$ws.Range("I137").Value2 = 1626.6
$ws.Range("K137").Value2 = 4879.799999999999
$ws.Range("M137").Value2 = -2329.799999999999
$ws.Range("H141").Value2 = 2160.12
$ws.Range("I141").Value2 = 1833.4584
$ws.Range("K141").Value2 = 5500.3752
$ws.Range("M141").Value2 = -320.3752000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 39817.035
$ws.Range("I74").Value2 = 78732.69500000001
$ws.Range("K74").Value2 = 78732.69500000001
$ws.Range("M74").Value2 = -77858.69500000001
$ws.Range("H77").Value2 = 39817.035
$ws.Range("I77").Value2 = 78732.69500000001
$ws.Range("K77").Value2 = 393663.475
$ws.Range("M77").Value2 = -389295.475
$ws.Range("H110").Value2 = 1999.1666
$ws.Range("I110").Value2 = 1665.3334
$ws.Range("K110").Value2 = 1665.3334
$ws.Range("M110").Value2 = 379.6666
$ws.Range("H132").Value2 = 2717.3684
$ws.Range("I132").Value2 = 2785.3572
$ws.Range("K132").Value2 = 8356.071599999999
$ws.Range("M132").Value2 = -5826.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value2 = 0
$ws.Range("I46").Value2 = 0
$ws.Range("K46").Value2 = 0
$ws.Range("M46").ClearContents()
$ws.Range("H86").Value2 = 2257.5
$ws.Range("I86").Value2 = 2224.889
$ws.Range("J86").Value2 = 2355.3333
$ws.Range("K86").Value2 = 2224.889
$ws.Range("L86").Value2 = 2355.3333
$ws.Range("M86").Value2 = -1101.889
$ws.Range("N86").Value2 = -4601.3333
$ws.Range("H89").Value2 = 2257.5
$ws.Range("I89").Value2 = 2224.889
$ws.Range("J89").Value2 = 2355.3333
$ws.Range("K89").Value2 = 11124.445
$ws.Range("L89").Value2 = 11776.6665
$ws.Range("M89").Value2 = -5508.445
$ws.Range("N89").Value2 = -23008.6665
$ws.Range("H107").Value2 = 2137.04
$ws.Range("I107").Value2 = 1974.1
$ws.Range("K107").Value2 = 1974.1
$ws.Range("M107").Value2 = -54.09999999999991
$ws.Range("H134").Value2 = 5037.7617
$ws.Range("J134").Value2 = 9040.5
$ws.Range("L134").Value2 = 27121.5
$ws.Range("N134").Value2 = -32191.5
$ws.Range("H140").Value2 = 76246
$ws.Range("J140").Value2 = 76246
$ws.Range("L140").Value2 = 76246
$ws.Range("N140").Value2 = -86606

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3326.0513
$ws.Range("I31").Value2 = 2774.52
$ws.Range("K31").Value2 = 2774.52
$ws.Range("M31").Value2 = -2479.52
$ws.Range("H34").Value2 = 3326.0513
$ws.Range("I34").Value2 = 2774.52
$ws.Range("K34").Value2 = 2774.52
$ws.Range("M34").Value2 = -2572.52
$ws.Range("H58").Value2 = 2117.5386
$ws.Range("I58").Value2 = 2026.4
$ws.Range("J58").Value2 = 2421.3333
$ws.Range("K58").Value2 = 2026.4
$ws.Range("L58").Value2 = 2421.3333
$ws.Range("M58").Value2 = -1823.4
$ws.Range("N58").Value2 = -2827.3333
$ws.Range("H136").Value2 = 2117.5386
$ws.Range("I136").Value2 = 2026.4
$ws.Range("J136").Value2 = 2421.3333
$ws.Range("K136").Value2 = 6079.200000000001
$ws.Range("L136").Value2 = 7263.999899999999
$ws.Range("M136").Value2 = -3529.200000000001
$ws.Range("N136").Value2 = -12363.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 3039700.5
$ws.Range("J113").Value2 = 4052435.8
$ws.Range("L113").Value2 = 12157307.4
$ws.Range("N113").Value2 = -12161647.4
$ws.Range("H132").Value2 = 7308.091
$ws.Range("J132").Value2 = 7308.091
$ws.Range("L132").Value2 = 65772.819
$ws.Range("N132").Value2 = -70832.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value2 = 15000.167
$ws.Range("J47").Value2 = 15000.167
$ws.Range("L47").Value2 = 15000.167
$ws.Range("N47").Value2 = -16136.167
$ws.Range("H107").Value2 = 929.1
$ws.Range("J107").Value2 = 901.6667
$ws.Range("L107").Value2 = 901.6667
$ws.Range("N107").Value2 = -4741.6667
$ws.Range("H132").Value2 = 6109.5186
$ws.Range("I132").Value2 = 5094
$ws.Range("K132").Value2 = 15282
$ws.Range("M132").Value2 = -12752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 1161.875
$ws.Range("I46").Value2 = 1109.6
$ws.Range("J46").Value2 = 1249
$ws.Range("K46").Value2 = 1109.6
$ws.Range("L46").Value2 = 1249
$ws.Range("M46").Value2 = -921.5999999999999
$ws.Range("N46").Value2 = -1625
$ws.Range("H61").Value2 = 0
$ws.Range("I61").Value2 = 0
$ws.Range("K61").Value2 = 0
$ws.Range("M61").ClearContents()
$ws.Range("H93").Value2 = 2435.1428
$ws.Range("I93").Value2 = 2174.6667
$ws.Range("K93").Value2 = 2174.6667
$ws.Range("M93").Value2 = -926.6667000000002
$ws.Range("H113").Value2 = 0
$ws.Range("I113").Value2 = 0
$ws.Range("K113").Value2 = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value2 = 2793.4
$ws.Range("I132").Value2 = 2847
$ws.Range("K132").Value2 = 8541
$ws.Range("M132").Value2 = -6011
$ws.Range("H136").Value2 = 2446.15
$ws.Range("I136").Value2 = 2090.0908
$ws.Range("J136").Value2 = 2881.3333
$ws.Range("K136").Value2 = 6270.2724
$ws.Range("L136").Value2 = 8643.999899999999
$ws.Range("M136").Value2 = -3720.2724
$ws.Range("N136").Value2 = -13743.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value2 = 43997.5
$ws.Range("I52").Value2 = 43997.5
$ws.Range("K52").Value2 = 43997.5
$ws.Range("M52").Value2 = -43771.5
$ws.Range("H58").Value2 = 40095.75
$ws.Range("I58").Value2 = 9997.5
$ws.Range("K58").Value2 = 9997.5
$ws.Range("M58").Value2 = -9689.5
$ws.Range("H61").Value2 = 1565342.8
$ws.Range("I61").Value2 = 1691616.6
$ws.Range("K61").Value2 = 1691616.6
$ws.Range("M61").Value2 = -1691324.6
$ws.Range("H113").Value2 = 2377.1
$ws.Range("I113").Value2 = 3179
$ws.Range("K113").Value2 = 9537
$ws.Range("M113").Value2 = -7367
$ws.Range("H122").Value2 = 1608.579
$ws.Range("I122").Value2 = 1401.2727
$ws.Range("K122").Value2 = 4203.8181
$ws.Range("M122").Value2 = -1753.8181
$ws.Range("H132").Value2 = 2627.3044
$ws.Range("I132").Value2 = 2557.6667
$ws.Range("K132").Value2 = 7673.000100000001
$ws.Range("M132").Value2 = -5143.000100000001
$ws.Range("H136").Value2 = 1214.5758
$ws.Range("I136").Value2 = 1056.2413
$ws.Range("J136").Value2 = 2362.5
$ws.Range("K136").Value2 = 3168.7239
$ws.Range("L136").Value2 = 7087.5
$ws.Range("M136").Value2 = -618.7239
$ws.Range("N136").Value2 = -12187.5
